# İş Takip Güncellemesi - 2.01.2026 15:35:42
# Shifts the tracked "İşe Başlama/Yer Teslimi", "İhale Bitiş Tarihi" and related
# "Güncelleme" planning dates back by one day, and appends the latest progress
# notes for the two rows (116 ve 118) still being worked on.

$wb = $excel.ActiveWorkbook
$wsList = $wb.Worksheets.Item(1)   # "İş Takip Listesi"
$wsUpd  = $wb.Worksheets.Item(2)   # "Güncelleme"

# --- "İş Takip Listesi" sheet: İŞE BAŞLAMA/YER TESLİMİ (J) & İHALE BİTİŞ TARİHİ (K) ---
$listDates = @(
    @("J2", "2025-06-19"),
    @("K2", "2025-11-19"),
    @("J3", "2025-06-19"),
    @("K3", "2025-11-19"),
    @("J4", "2025-06-19"),
    @("K4", "2025-11-19"),
    @("J5", "2025-06-19"),
    @("K5", "2025-11-19"),
    @("J6", "2025-06-19"),
    @("K6", "2025-11-19"),
    @("J7", "2025-06-19"),
    @("K7", "2025-11-19"),
    @("J8", "2025-06-19"),
    @("K8", "2025-11-19"),
    @("J9", "2025-06-19"),
    @("K9", "2025-11-19"),
    @("J10", "2025-06-19"),
    @("K10", "2025-11-19"),
    @("J33", "2025-06-21"),
    @("K33", "2025-11-21"),
    @("J34", "2025-06-21"),
    @("K34", "2025-11-21"),
    @("J35", "2025-06-21"),
    @("K35", "2025-11-21"),
    @("J36", "2025-06-21"),
    @("K36", "2025-11-21"),
    @("J37", "2025-06-21"),
    @("K37", "2025-11-21"),
    @("J38", "2025-06-21"),
    @("K38", "2025-11-21"),
    @("J39", "2025-06-21"),
    @("K39", "2025-11-21"),
    @("J40", "2025-06-21"),
    @("K40", "2025-11-21"),
    @("J41", "2025-06-21"),
    @("K41", "2025-11-21"),
    @("J42", "2025-06-21"),
    @("K42", "2025-11-21"),
    @("J43", "2025-06-21"),
    @("K43", "2025-11-21"),
    @("J44", "2025-06-21"),
    @("K44", "2025-11-21"),
    @("J45", "2025-06-21"),
    @("K45", "2025-11-21"),
    @("J46", "2025-06-21"),
    @("K46", "2025-11-21"),
    @("J47", "2025-06-21"),
    @("K47", "2025-11-21"),
    @("J48", "2025-06-21"),
    @("K48", "2025-11-21"),
    @("J49", "2025-06-21"),
    @("K49", "2025-11-21"),
    @("J50", "2025-06-21"),
    @("K50", "2025-11-21"),
    @("J51", "2025-06-21"),
    @("K51", "2025-11-21"),
    @("J52", "2025-06-21"),
    @("K52", "2025-11-21"),
    @("J53", "2025-06-21"),
    @("K53", "2025-11-21"),
    @("J54", "2025-06-21"),
    @("K54", "2025-11-21"),
    @("J55", "2025-06-21"),
    @("K55", "2025-11-21"),
    @("J56", "2025-06-21"),
    @("K56", "2025-11-21"),
    @("J57", "2025-06-21"),
    @("K57", "2025-11-21"),
    @("J58", "2025-06-21"),
    @("K58", "2025-11-21"),
    @("J59", "2025-06-21"),
    @("K59", "2025-11-21"),
    @("J60", "2025-06-21"),
    @("K60", "2025-11-21"),
    @("J61", "2025-06-21"),
    @("K61", "2025-11-21"),
    @("J62", "2025-06-21"),
    @("K62", "2025-11-21"),
    @("J63", "2025-06-21"),
    @("K63", "2025-11-21"),
    @("J64", "2025-06-21"),
    @("K64", "2025-11-21"),
    @("J65", "2025-06-21"),
    @("K65", "2025-11-21"),
    @("J66", "2025-06-21"),
    @("K66", "2025-11-21"),
    @("J67", "2025-06-21"),
    @("K67", "2025-11-21"),
    @("J68", "2025-06-21"),
    @("K68", "2025-11-21"),
    @("J69", "2025-06-21"),
    @("K69", "2025-11-21"),
    @("J70", "2025-06-21"),
    @("K70", "2025-11-21"),
    @("J71", "2025-06-21"),
    @("K71", "2025-11-21"),
    @("J72", "2025-06-21"),
    @("K72", "2025-11-21"),
    @("J73", "2025-06-21"),
    @("K73", "2025-11-21"),
    @("J74", "2025-06-21"),
    @("K74", "2025-11-21"),
    @("J75", "2025-06-21"),
    @("K75", "2025-11-21"),
    @("J76", "2025-06-21"),
    @("K76", "2025-11-21"),
    @("J77", "2025-06-21"),
    @("K77", "2025-11-21"),
    @("J78", "2025-06-21"),
    @("K78", "2025-11-21"),
    @("J79", "2025-06-21"),
    @("K79", "2025-11-21"),
    @("J80", "2025-06-21"),
    @("K80", "2025-11-21"),
    @("J81", "2025-06-21"),
    @("K81", "2025-11-21"),
    @("J82", "2025-06-21"),
    @("K82", "2025-11-21"),
    @("J83", "2025-06-21"),
    @("K83", "2025-11-21"),
    @("J84", "2025-06-21"),
    @("K84", "2025-11-21"),
    @("J85", "2025-06-21"),
    @("K85", "2025-11-21"),
    @("J86", "2025-06-21"),
    @("K86", "2025-11-21"),
    @("J87", "2025-06-21"),
    @("K87", "2025-11-21"),
    @("J88", "2025-06-21"),
    @("K88", "2025-11-21"),
    @("J89", "2025-06-21"),
    @("K89", "2025-11-21"),
    @("J90", "2025-06-21"),
    @("K90", "2025-11-21"),
    @("J91", "2025-06-21"),
    @("K91", "2025-11-21"),
    @("J92", "2025-06-21"),
    @("K92", "2025-11-21"),
    @("J93", "2025-06-21"),
    @("K93", "2025-11-21"),
    @("J94", "2025-06-21"),
    @("K94", "2025-11-21"),
    @("J95", "2024-04-19"),
    @("K95", "2025-06-13"),
    @("J96", "2024-04-19"),
    @("K96", "2025-06-13"),
    @("J97", "2024-04-19"),
    @("K97", "2025-06-13"),
    @("J98", "2024-04-19"),
    @("K98", "2025-06-13"),
    @("J99", "2024-04-19"),
    @("K99", "2025-06-13"),
    @("J100", "2024-04-19"),
    @("K100", "2025-06-13"),
    @("J101", "2024-04-19"),
    @("K101", "2025-06-13"),
    @("J102", "2024-04-19"),
    @("K102", "2025-06-13"),
    @("J103", "2024-04-19"),
    @("K103", "2025-06-13"),
    @("J104", "2024-04-19"),
    @("K104", "2025-06-13"),
    @("J105", "2024-04-19"),
    @("K105", "2025-06-13"),
    @("J106", "2024-04-19"),
    @("K106", "2025-06-13"),
    @("J107", "2024-04-19"),
    @("K107", "2025-06-13"),
    @("J108", "2024-04-19"),
    @("K108", "2025-06-13"),
    @("J109", "2024-04-19"),
    @("K109", "2025-06-13"),
    @("J110", "2024-04-19"),
    @("K110", "2025-06-13"),
    @("J111", "2024-04-19"),
    @("K111", "2025-06-13"),
    @("J112", "2024-04-19"),
    @("K112", "2025-06-13"),
    @("J113", "2024-04-19"),
    @("K113", "2025-06-13"),
    @("J114", "2024-04-19"),
    @("K114", "2025-06-13"),
    @("J115", "2024-04-19"),
    @("K115", "2025-06-13"),
    @("J116", "2024-04-19"),
    @("K116", "2025-06-13"),
    @("J117", "2024-04-19"),
    @("K117", "2025-06-13"),
    @("J118", "2024-04-19"),
    @("K118", "2025-06-13"),
    @("J119", "2024-04-19"),
    @("K119", "2025-06-13"),
    @("J120", "2024-04-19"),
    @("K120", "2025-06-13"),
    @("J121", "2024-04-19"),
    @("K121", "2025-06-13"),
    @("J122", "2024-04-19"),
    @("K122", "2025-06-13")
)
foreach ($pair in $listDates) {
    $cell = $wsList.Range($pair[0])
    $cell.NumberFormat = "@"   # keep the literal yyyy-mm-dd text, not a real date serial
    $cell.Value = $pair[1]
}

# --- "Güncelleme" sheet: UÇUŞ TARİHİ (I/J), BİLGİLENDİRME İLANI (N) & KESİN ASKI (P) ---
$updDates = @(
    @("J2", "2024-07-26"),
    @("N2", "2025-03-29"),
    @("P2", "2025-06-15"),
    @("J3", "2024-10-27"),
    @("N3", "2025-07-17"),
    @("P3", "2025-10-25"),
    @("J4", "2024-08-31"),
    @("N4", "2025-02-21"),
    @("P4", "2025-05-18"),
    @("I5", "2025-02-24"),
    @("J6", "2025-10-07"),
    @("N6", "2025-06-26"),
    @("P6", "2025-12-22"),
    @("I7", "2024-10-27"),
    @("J7", "2024-10-27"),
    @("J8", "2024-10-15"),
    @("N8", "2025-03-16"),
    @("P8", "2025-04-18"),
    @("I9", "2025-06-12"),
    @("J9", "2024-11-28"),
    @("J10", "2024-09-26"),
    @("N10", "2025-07-07"),
    @("P10", "2025-10-29"),
    @("I11", "2025-04-03"),
    @("J11", "2024-11-09"),
    @("N11", "2025-07-27"),
    @("P11", "2025-12-22"),
    @("J12", "2024-10-07"),
    @("N12", "2025-06-16"),
    @("P12", "2025-10-19"),
    @("J13", "2024-12-05"),
    @("J14", "2025-10-03"),
    @("N14", "2025-10-25"),
    @("J15", "2024-12-24"),
    @("N15", "2025-07-14"),
    @("P15", "2025-10-26"),
    @("J16", "2024-08-22"),
    @("N16", "2025-01-30"),
    @("P16", "2025-04-18"),
    @("J17", "2024-09-07"),
    @("N17", "2025-10-25"),
    @("J18", "2025-02-14"),
    @("I19", "2025-04-04"),
    @("J19", "2024-12-24"),
    @("N19", "2025-08-03"),
    @("J20", "2024-12-05"),
    @("N20", "2025-12-09"),
    @("J21", "2024-09-28"),
    @("J22", "2024-09-28"),
    @("J23", "2024-12-06"),
    @("I24", "2025-06-02"),
    @("J25", "2024-10-31"),
    @("J27", "2025-01-21"),
    @("J28", "2024-11-19"),
    @("N28", "2025-11-04"),
    @("I29", "2025-02-08"),
    @("J29", "2024-12-06"),
    @("N29", "2025-10-16")
)
foreach ($pair in $updDates) {
    $cell = $wsUpd.Range($pair[0])
    $cell.NumberFormat = "@"   # keep the literal yyyy-mm-dd text, not a real date serial
    $cell.Value = $pair[1]
}

# --- "İş Takip Listesi" sheet: append latest progress notes (NOTLAR column) ---
$wsList.Range("M116").Value = "01.11.2025 firmaya teslim edilecek`n05.11.2025 kontroller devam ediyor`n13.11.2025 Firmaya 28.11.2025 teslim edilecek`n02.12.2025 Firmaya 05.12.2025 teslim edilecek`n15.12.2025 Değerlendirme devam ediyor`n02.01.2026 Firmaya 15.01.2026 Teslim edilecek"
$wsList.Range("M118").Value = "01.11.2025 firmaya teslim edilecek`n05.11.2025 kontroller devam ediyor`n10.11.2025 Değerlendrime 13.11.2025 de bitecek`n13.11.2025 Firmaya 17-21 haftası teslim edilecek`n20.11.2025 firmaya 24.11.2025 de teslim edilecek`n02.12.2025 bugün firmaya teslim edilecek`n15.12.2025 Proje hazırlanıyor`n22.12.2025 Proje 25.12.2025 biticek, krokiler 30.12.2025 verilecek`n02.01.2026 Firma projeyi hazırlıyor "

# Re-fit the two note rows so their height tracks the (now one-line-longer) text,
# same as Excel does automatically when a wrapped cell grows.
$wsList.Rows.Item(116).AutoFit()
$wsList.Rows.Item(118).AutoFit()

Write-Output "Is Takip guncellemesi uygulandi."
